# Apply latest code change for github run
# - Admin sheet: regenerate the test username value in D2
# - Jira sheet: update the automation run log (grown from 4 to 6 rows),
#   clearing the stale "Error description" entries for the Key /
#   PersonalDetails rows and appending the newly-failed test keys.

$wb = $excel.ActiveWorkbook

# --- Admin sheet: regenerate test username value in D2 ---
$admin = $wb.Worksheets.Item("Admin")
$admin.Range("D2").Value = "7598587A"

# --- Jira sheet: update run log rows ---
$jira = $wb.Worksheets.Item("Jira")

# Row 2 ("Key"): clear the stale error description, keeping the cell as an
# (empty) text value rather than wiping it out entirely.
$b2 = $jira.Cells.Item(2, 2)
$b2.Value = "'"
$b2.Style = "Normal"

# Row 3 now references PersonalDetails, also with a blank description.
$jira.Cells.Item(3, 1).Value = "PersonalDetails"
$b3 = $jira.Cells.Item(3, 2)
$b3.Value = "'"
$b3.Style = "Normal"

# Row 4: Recruitment_HiredList
$jira.Cells.Item(4, 1).Value = "Recruitment_HiredList"
$jira.Cells.Item(4, 2).Style = "Normal"

# New row 5: Recruitment_RejectionList1
$jira.Cells.Item(5, 1).Value = "Recruitment_RejectionList1"
$jira.Cells.Item(5, 2).Style = "Normal"

# New row 6: PersonalDetails
$jira.Cells.Item(6, 1).Value = "PersonalDetails"
$jira.Cells.Item(6, 2).Style = "Normal"
